# Weekly price update: insert a new report row for "Cardinal" potato
# (1a cosecha, Provincia del Elquí) at the top of the data block, at
# row 646, pushing the existing rows 646-677 down to 647-678.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 646 - shifts rows 646:677 down to 647:678
# and extends the used range to A1:R678.
$ws.Rows.Item(646).Insert()

# Populate the new row 646 with the latest weekly report values.
$ws.Range("A646").Value = 8
$ws.Range("B646").Value = "Terminal La Palmera de La Serena"
$ws.Range("C646").Value = "Coquimbo"
$ws.Range("D646").Value = 45119
$ws.Range("E646").Value = 4
$ws.Range("F646").Value = 100114001
$ws.Range("G646").Value = "Papa"
$ws.Range("H646").Value = "Cardinal"
$ws.Range("I646").Value = "1a (cosecha)"
$ws.Range("J646").Value = 600
$ws.Range("K646").Value = 19000
$ws.Range("L646").Value = 19500
$ws.Range("M646").Value = 19250
$ws.Range("N646").Value = "`$/saco 25 kilos"
$ws.Range("O646").Value = "Provincia del Elquí"
$ws.Range("P646").Value = 770
$ws.Range("Q646").Value = 25
$ws.Range("R646").Value = "Hortaliza"
